function Make-Run($text) {
    if ($text -ne $text.Trim()) {
        return '<w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    } else {
        return '<w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>' + $text + '</w:t></w:r>'
    }
}

$d = $word.ActiveDocument

# Locate the target paragraph - the one that starts the "Primero debo esperar..."
# reflection text near the end of the document.
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Primero debo esperar a que el infiltrado")) {
        $para = $p
        break
    }
}

if ($para -eq $null) {
    throw "Target paragraph not found"
}

# Keep the paragraph's own opening-tag attributes (paraId/textId/rsid*) so we
# preserve its identity while rewriting its run content and appending two
# brand-new paragraphs right after it (mirroring how Word keeps the original
# paragraph mark for edited text and leaves freshly-typed paragraphs unstamped).
$pPrefix = '<w:p w14:paraId="4D1BABDE" w14:textId="5EE56383" w:rsidR="00850BFC" w:rsidRDefault="00850BFC" w:rsidP="00515072" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">'
$pPr = '<w:pPr><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr>'

$xml = $pPrefix
$xml += $pPr
$xml += Make-Run "Primero debo esperar a que el infiltrado me diga los parámetros de salida del cañón O para ver si en realidad es una amenaza, ya que se puede dar que con los parámetros de salida del cañón O, la bala no tenga la suficiente fuerza para impactar"
$xml += Make-Run " o estar a rango"
$xml += Make-Run " "
$xml += Make-Run "de"
$xml += Make-Run " dañar mi cañón defensivo"
$xml += Make-Run " y para no gastar una bala innecesariamente no habría necesidad de disparar."
$xml += '</w:p>'

$xml += '<w:p>'
$xml += $pPr
$xml += Make-Run "Si con los parámetros de salida del cañón O, analizo que verdaderamente puede ser una amenaza, debo determinar "
$xml += Make-Run "en qué lugar podría impactar la bala, luego de tener el valor, darle la indicación a mi cañón defensivo que prepare una bala para detener la bala del cañón O sin que el impacto dañe la integridad de ambos cañones."
$xml += '</w:p>'

$xml += '<w:p>'
$xml += $pPr
$xml += Make-Run "Debo determinar que el rango de explosión no dañe ninguno de los dos cañones."
$xml += '</w:p>'

# Replace the whole paragraph (including its end-of-paragraph mark) with the
# rebuilt paragraph plus the two new ones that follow it.
$full = $d.Range($para.Range.Start, $para.Range.End)
$full.InsertXML($xml)
